# Checks and fixes of charging strategies with many trucks and limited grid cap
# - adjusted input files (no fixed charging profile, more vehicles and smaller grid connection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete 400km/500km trip-pattern rows (delete bottom-up so the
# remaining row numbers don't shift while we work):
#   row 12 -> 6_18.5_500km
#   row 11 -> 6_18.5_400km
#   row 9  -> 6_17_500km
#   row 8  -> 6_17_400km
#   row 5  -> 6_16_500km
#   row 4  -> 6_16_400km
$ws.Rows(12).Delete()
$ws.Rows(11).Delete()
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
$ws.Rows(4).Delete()

# After the deletions, the sheet has 5 trip patterns left (rows 2-6):
#   6_16_250km, 6_16_300km, 6_17_250km, 6_17_300km, 6_18.5_300km
# Append the 3 new trip patterns as rows 7-9.

$ws.Range("A7").Value = "6.5_18.25_280km"
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 360
$ws.Range("D7").Value = 1080
$ws.Range("E7").Value = 291
$ws.Range("F7").Value = 1800
$ws.Range("G7").Value = 2520
$ws.Range("H7").Value = 282
$ws.Range("I7").Value = 3300
$ws.Range("J7").Value = 3960
$ws.Range("K7").Value = 285
$ws.Range("L7").Value = 4680
$ws.Range("M7").Value = 5400
$ws.Range("N7").Value = 265
$ws.Range("O7").Value = 6120
$ws.Range("P7").Value = 6840
$ws.Range("Q7").Value = 301
$ws.Range("R7").Value = 7560
$ws.Range("S7").Value = 8280
$ws.Range("T7").Value = 268
$ws.Range("U7").Value = 9000
$ws.Range("V7").Value = 9720
$ws.Range("W7").Value = 260

$ws.Range("A8").Value = "6_18_250km"
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 360
$ws.Range("D8").Value = 1080
$ws.Range("E8").Value = 230
$ws.Range("F8").Value = 1800
$ws.Range("G8").Value = 2520
$ws.Range("H8").Value = 235
$ws.Range("I8").Value = 3240
$ws.Range("J8").Value = 3960
$ws.Range("K8").Value = 241
$ws.Range("L8").Value = 4680
$ws.Range("M8").Value = 5400
$ws.Range("N8").Value = 276
$ws.Range("O8").Value = 6120
$ws.Range("P8").Value = 6840
$ws.Range("Q8").Value = 250
$ws.Range("R8").Value = 7560
$ws.Range("S8").Value = 8280
$ws.Range("T8").Value = 240

$ws.Range("A9").Value = "6.5_19.5_340km"
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 420
$ws.Range("D9").Value = 1140
$ws.Range("E9").Value = 330
$ws.Range("F9").Value = 1800
$ws.Range("G9").Value = 2640
$ws.Range("H9").Value = 323
$ws.Range("I9").Value = 3300
$ws.Range("J9").Value = 4080
$ws.Range("K9").Value = 355
$ws.Range("L9").Value = 4740
$ws.Range("M9").Value = 5460
$ws.Range("N9").Value = 335
$ws.Range("O9").Value = 6120
$ws.Range("P9").Value = 6960
$ws.Range("Q9").Value = 327
$ws.Range("R9").Value = 7620
$ws.Range("S9").Value = 8400
$ws.Range("T9").Value = 328
$ws.Range("U9").Value = 9000
$ws.Range("V9").Value = 9840
$ws.Range("W9").Value = 350

# Update the selection shown when the file is reopened
$ws.Range("E14").Select()
